$d = $word.ActiveDocument

# 1. Remove the "Das Mietverhältnis beginnt am 01.04.2023. " prefix from the paragraph
$d.Content.Find.Execute("Das Mietverhältnis beginnt am 01.04.2023. Das Mietverhältnis läuft auf unbestimmte Zeit.", $true, $false, $false, $false, $false, $true, 1, $false, "Das Mietverhältnis läuft auf unbestimmte Zeit.", 2)

# 2. Remove the word "solcher " before "Befristungsgrund"
$d.Content.Find.Execute("Besteht kein solcher Befristungsgrund", $true, $false, $false, $false, $false, $true, 1, $false, "Besteht kein Befristungsgrund", 2)
